$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("watchlist")
$ws2 = $wb.Worksheets.Item("stocks")
$ws2.Range("F1:F15").Delete(-4159)
$ws2.Activate()
$ws2.Range("H1").Select()
